$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 8: Loading detail values (string replacements) ---
# "FB800-1,IOB800-1,POS800-S-1,POS800-M-1" combo strings were split into single values
$ws.Range("G8").Value = "PCH800 5.0A"
$ws.Range("H8").Value = "PCH800 5.0A-1"
$ws.Range("L8").Value = "Fuse board - 1,IOB800 - 1,POS800-S - 1,POS800-M - 1"

# --- Row 7: two new header columns ---
$hdrSrc = $ws.Range("N7:O7")
$hdrDst = $ws.Range("Q7:R7")
$hdrSrc.Copy()
$hdrDst.PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q7").Value = "Expected Max 24V PSU Load value"
$ws.Range("R7").Value = "Expected Max 2nd 24V PSU Load value"

# --- Row 8: two new data columns (quote-prefixed numeric-looking text) ---
$cellSrc = $ws.Range("C5")
$cellDst = $ws.Range("Q8:R8")
$cellSrc.Copy()
$cellDst.PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q8").Value = "'9.700"
$ws.Range("R8").Value = "'5.000"

# --- Row 8: updated numeric load values ---
$ws.Range("F8").Value = 0.319
$ws.Range("J8").Value = 0.319
$ws.Range("N8").Value = 1.023
$ws.Range("O8").Value = 1.023

$excel.CutCopyMode = 0

# --- Column widths for the two new columns ---
$ws.Columns.Item(17).ColumnWidth = 29
$ws.Columns.Item(18).ColumnWidth = 30

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 80
$ws.Range("P8").Select()
